# Updates cryptocurrency price/volume figures on the Cryptos sheet,
# matching the GitHub Actions scrape commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.062.19'
$ws.Range("E2").Value = '  +2.28%  '
$ws.Range("D3").Value = '2.302.00'
$ws.Range("E3").Value = '  +1.61%  '
$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = '''301.67'
$ws.Range("E5").Value = '  +1.25%  '
$ws.Range("D6").Value = '''99.38'
$ws.Range("E6").Value = '  +5.56%  '
$ws.Range("D7").Value = '''0.507'
$ws.Range("E7").Value = '  +1.99%  '
$ws.Range("E8").Value = '  -0.08%  '
$ws.Range("E9").Value = '  +3.10%  '
$ws.Range("D10").Value = '''34.24'
$ws.Range("E10").Value = '  +3.82%  '
$ws.Range("D11").Value = '''0.0799'
$ws.Range("E11").Value = '  +1.14%  '
$ws.Range("D12").Value = '''49.05'
$ws.Range("E12").Value = '  +1.74%  '
$ws.Range("E13").Value = '  +4.30%  '
$ws.Range("D14").Value = '''17.96'
$ws.Range("E14").Value = '  +15.56%  '
$ws.Range("D15").Value = '''6.80'
$ws.Range("E15").Value = '  +2.17%  '
$ws.Range("D16").Value = '2.657.60'
$ws.Range("E16").Value = '  +1.58%  '
$ws.Range("D17").Value = '2.281.27'
$ws.Range("E17").Value = '  +1.32%  '
$ws.Range("E18").Value = '  +4.69%  '
$ws.Range("D19").Value = '42.968.81'
$ws.Range("E19").Value = '  +2.07%  '
$ws.Range("D20").Value = '''12.52'
$ws.Range("E20").Value = '  +10.48%  '
$ws.Range("D21").Value = '0.0₃0907'
$ws.Range("E21").Value = '  +2.01%  '
$ws.Range("E22").Value = '  +1.60%  '
$ws.Range("D23").Value = '''67.71'
$ws.Range("E23").Value = '  +1.55%  '
$ws.Range("D24").Value = '''236.57'
$ws.Range("E24").Value = '  +1.45%  '
$ws.Range("E25").Value = '  +14.69%  '
$ws.Range("E26").Value = '  +0.09%  '
$ws.Range("E27").Value = '  +0.51%  '
$ws.Range("D28").Value = '''24.75'
$ws.Range("E28").Value = '  +4.02%  '
$ws.Range("D29").Value = '''167.87'
$ws.Range("E29").Value = '  +0.12%  '
$ws.Range("E30").Value = '  -9.12%  '
$ws.Range("D31").Value = '''34.00'
$ws.Range("E31").Value = '  +1.19%  '
$ws.Range("D32").Value = '''9.14'
$ws.Range("E32").Value = '  +1.29%  '
$ws.Range("D33").Value = '''0.999'
$ws.Range("E33").Value = '  -0.04%  '
$ws.Range("D34").Value = '''5.05'
$ws.Range("E34").Value = '  +2.56%  '
$ws.Range("D35").Value = '''2.44'
$ws.Range("E35").Value = '  +4.36%  '
$ws.Range("D36").Value = '''4.56'
$ws.Range("E36").Value = '  +1.58%  '
$ws.Range("D37").Value = '''17.01'
$ws.Range("E37").Value = '  +5.09%  '
$ws.Range("D38").Value = '''0.0692'
$ws.Range("E38").Value = '  +0.38%  '
$ws.Range("E39").Value = '  +3.93%  '
$ws.Range("E40").Value = '  +5.19%  '
$ws.Range("D41").Value = '''2.80'
$ws.Range("E41").Value = '  +1.22%  '
$ws.Range("D42").Value = '''0.110'
$ws.Range("E42").Value = '  +0.29%  '
$ws.Range("E43").Value = '  -5.33%  '
$ws.Range("D44").Value = '1.993.06'
$ws.Range("E44").Value = '  +1.69%  '
$ws.Range("D45").Value = '''0.0285'
$ws.Range("E45").Value = '  +2.72%  '
$ws.Range("D46").Value = '''10.04'
$ws.Range("E46").Value = '  +5.20%  '
$ws.Range("D47").Value = '''17.73'
$ws.Range("E47").Value = '  +2.31%  '
$ws.Range("E48").Value = '  +2.93%  '
$ws.Range("D49").Value = '''56.37'
$ws.Range("E49").Value = '  +8.78%  '
$ws.Range("D50").Value = '2.526.59'
$ws.Range("E50").Value = '  +1.46%  '
$ws.Range("D51").Value = '''1.54'
$ws.Range("E51").Value = '  +4.51%  '
